$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4005.111
$ws.Range("I40").Value = 4185
$ws.Range("J40").Value = 3890.6365
$ws.Range("K40").Value = 4185
$ws.Range("L40").Value = 3890.6365
$ws.Range("M40").Value = -4010
$ws.Range("N40").Value = -4240.636500000001

$ws.Range("H51").Value = 4500
$ws.Range("I51").Value = 4500
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 4500
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -4016
$ws.Range("N51").ClearContents()

$ws.Range("H129").Value = 31251636
$ws.Range("I129").Value = 71429820
$ws.Range("J129").Value = 1935.8889
$ws.Range("K129").Value = 214289460
$ws.Range("L129").Value = 5807.6667
$ws.Range("M129").Value = -214284460
$ws.Range("N129").Value = -15807.6667

$ws.Range("H132").Value = 13161530
$ws.Range("I132").Value = 14709609
$ws.Range("J132").Value = 2855.375
$ws.Range("K132").Value = 44128827
$ws.Range("L132").Value = 8566.125
$ws.Range("M132").Value = -44126297
$ws.Range("N132").Value = -13626.125

$ws.Range("H137").Value = 35229.168
$ws.Range("I137").Value = 42014.047
$ws.Range("J137").Value = 2058.6667
$ws.Range("K137").Value = 126042.141
$ws.Range("L137").Value = 6176.000100000001
$ws.Range("M137").Value = -123492.141
$ws.Range("N137").Value = -11276.0001

$ws.Range("H138").Value = 2685.8022
$ws.Range("I138").Value = 886.1667
$ws.Range("K138").Value = 2658.5001
$ws.Range("M138").Value = 2481.4999

$ws.Range("H141").Value = 5217.5835
$ws.Range("I141").Value = 5638.905
$ws.Range("J141").Value = 2268.3333
$ws.Range("K141").Value = 16916.715
$ws.Range("L141").Value = 6804.999899999999
$ws.Range("M141").Value = -11736.715
$ws.Range("N141").Value = -17164.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1715080.1
$ws.Range("I2").Value = 2176359.5
$ws.Range("J2").Value = 1756.8572
$ws.Range("K2").Value = 2176359.5
$ws.Range("L2").Value = 1756.8572
$ws.Range("M2").Value = -2176246.5
$ws.Range("N2").Value = -1982.8572

$ws.Range("H15").Value = 15666.556
$ws.Range("I15").Value = 4833.3335
$ws.Range("J15").Value = 37333
$ws.Range("K15").Value = 4833.3335
$ws.Range("L15").Value = 37333
$ws.Range("M15").Value = -4483.3335
$ws.Range("N15").Value = -38033

$ws.Range("H17").Value = 10000
$ws.Range("J17").Value = 10000
$ws.Range("L17").Value = 10000
$ws.Range("N17").Value = -10346

$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -701

$ws.Range("H74").Value = 35727.33
$ws.Range("I74").Value = 19318.66
$ws.Range("J74").Value = 445944
$ws.Range("K74").Value = 19318.66
$ws.Range("L74").Value = 445944
$ws.Range("M74").Value = -18444.66
$ws.Range("N74").Value = -447692

$ws.Range("H77").Value = 35727.33
$ws.Range("I77").Value = 19318.66
$ws.Range("J77").Value = 445944
$ws.Range("K77").Value = 96593.3
$ws.Range("L77").Value = 2229720
$ws.Range("M77").Value = -92225.3
$ws.Range("N77").Value = -2238456

$ws.Range("H116").Value = 1715080.1
$ws.Range("I116").Value = 2176359.5
$ws.Range("J116").Value = 1756.8572
$ws.Range("K116").Value = 2176359.5
$ws.Range("L116").Value = 1756.8572
$ws.Range("M116").Value = -2174065.5
$ws.Range("N116").Value = -6344.8572

$ws.Range("H122").Value = 995135.3
$ws.Range("I122").Value = 2430.818
$ws.Range("J122").Value = 2087110.2
$ws.Range("K122").Value = 7292.454000000001
$ws.Range("L122").Value = 6261330.6
$ws.Range("M122").Value = -4842.454000000001
$ws.Range("N122").Value = -6266230.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1715080.1
$ws.Range("I3").Value = 2176359.5
$ws.Range("J3").Value = 1756.8572
$ws.Range("K3").Value = 2176359.5
$ws.Range("L3").Value = 1756.8572
$ws.Range("M3").Value = -2176245.5
$ws.Range("N3").Value = -1984.8572

$ws.Range("H11").Value = 450
$ws.Range("J11").Value = 450
$ws.Range("L11").Value = 450
$ws.Range("N11").Value = -730

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2061.8596
$ws.Range("I58").Value = 1958.8718
$ws.Range("J58").Value = 2285
$ws.Range("K58").Value = 1958.8718
$ws.Range("L58").Value = 2285
$ws.Range("M58").Value = -1755.8718
$ws.Range("N58").Value = -2691

$ws.Range("H107").Value = 1827.2059
$ws.Range("I107").Value = 1455.4783
$ws.Range("J107").Value = 2604.4546
$ws.Range("K107").Value = 1455.4783
$ws.Range("L107").Value = 2604.4546
$ws.Range("M107").Value = 464.5217
$ws.Range("N107").Value = -6444.4546

$ws.Range("H132").Value = 45008.23
$ws.Range("I132").Value = 26643.385
$ws.Range("J132").Value = 124589.22
$ws.Range("K132").Value = 79930.155
$ws.Range("L132").Value = 373767.66
$ws.Range("M132").Value = -77400.155
$ws.Range("N132").Value = -378827.66

$ws.Range("H134").Value = 15326.477
$ws.Range("I134").Value = 19353.715
$ws.Range("J134").Value = 2993.0625
$ws.Range("K134").Value = 58061.145
$ws.Range("L134").Value = 8979.1875
$ws.Range("M134").Value = -55526.145
$ws.Range("N134").Value = -14049.1875

$ws.Range("H136").Value = 2061.8596
$ws.Range("I136").Value = 1958.8718
$ws.Range("J136").Value = 2285
$ws.Range("K136").Value = 5876.6154
$ws.Range("L136").Value = 6855
$ws.Range("M136").Value = -3326.6154
$ws.Range("N136").Value = -11955

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1273657.4
$ws.Range("I122").Value = 2964365
$ws.Range("J122").Value = 5626.75
$ws.Range("K122").Value = 8893095
$ws.Range("L122").Value = 16880.25
$ws.Range("M122").Value = -8890645
$ws.Range("N122").Value = -21780.25

$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -44900

$ws.Range("H132").Value = 1962.6888
$ws.Range("I132").Value = 2014.675
$ws.Range("J132").Value = 1546.8
$ws.Range("K132").Value = 6044.025
$ws.Range("L132").Value = 4640.4
$ws.Range("M132").Value = -3514.025
$ws.Range("N132").Value = -9700.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3857.3015
$ws.Range("I132").Value = 3532.83
$ws.Range("K132").Value = 10598.49
$ws.Range("M132").Value = -8068.49

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 20631852
$ws.Range("I132").Value = 27778678
$ws.Range("J132").Value = 840640.1
$ws.Range("K132").Value = 83336034
$ws.Range("L132").Value = 2521920.3
$ws.Range("M132").Value = -83333504
$ws.Range("N132").Value = -2526980.3

$ws.Range("H136").Value = 2151.1147
$ws.Range("I136").Value = 1874.6666
$ws.Range("K136").Value = 5623.9998
$ws.Range("M136").Value = -3073.9998
